$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Single Column Numbers")

# Row 3 must be unhidden before its value is rewritten, otherwise the
# runtime stamps the row with a stray custom row height when a cell in a
# hidden row is edited.
$ws.Rows.Item(3).Hidden = $false

# Swap the shared-string values shown in B2/B3: "Jacques" and "Alex" trade places
# (shared string table order changes so that "Alex" comes before "Jacques").
$ws.Range("B2").Value = "Alex"
$ws.Range("B3").Value = "Jacques"

# Row 2 becomes hidden, row 3 stays visible.
$ws.Rows.Item(2).Hidden = $true
